$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("A1").Value = "Ordem autorização"
$ws.Range("B1").Value = "ID"

# Update existing data row 2
$ws.Range("A2").Value = 1037083
$ws.Range("B2").Value = 34130918

# Add new data row 3
$ws.Range("A3").Value = 1036133
$ws.Range("B3").Value = 34131096

# Update selection to match target state
$ws.Range("B9").Select()
